$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.045860886573792
$ws.Range("B1").Value = 2.249685287475586
$ws.Range("C1").Value = 4.219629764556885
$ws.Range("D1").Value = 0.8181212544441223
$ws.Range("E1").Value = 1.093196034431458
